$wb = $excel.ActiveWorkbook

# Every "variables_*" sheet lists a handful of category names in column A.
# One category, "congenital", was renamed to "misc_long_term" across all
# the sheets that reference it.
$oldValue = "congenital"
$newValue = "misc_long_term"

$targetSheets = @(
    "variables_268", "variables_269", "variables_270", "variables_271",
    "variables_272", "variables_273", "variables_274",
    "variables_296", "variables_297", "variables_298", "variables_299",
    "variables_300", "variables_301", "variables_302", "variables_303",
    "variables_304", "variables_305", "variables_306", "variables_307",
    "variables_308", "variables_309", "variables_310", "variables_311",
    "variables_312", "variables_313", "variables_314", "variables_315",
    "variables_316",
    "variables_352", "variables_353", "variables_354", "variables_355",
    "variables_356"
)

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    for ($row = 1; $row -le $used.Rows.Count; $row++) {
        $cell = $ws.Cells.Item($row, 1)
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value = $newValue
            break
        }
    }
}
